$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "58.202.76"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.77%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.464.70"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "511.18"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.23%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "134.12"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.29%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.464.26"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0985"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  -5.73%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.901.71"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.82%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "58.060.63"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "22.06"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("E17").Value = "  +2.29%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.383.35"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  +0.75%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "316.13"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  +6.50%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "65.55"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("E29").Value = "  +5.89%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "171.89"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.76%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0740"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.70"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +0.12%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "18.12"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("E39").Value = "  +4.24%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.86"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  +2.80%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.816"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.05%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "136.84"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +12.48%  "
$ws.Range("E44").Value = "  +1.37%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.94"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.82%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.577"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "256.71"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0918"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0496"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.53%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.30"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.11%  "
